$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace course entry (EDUC_7094 -> SIST_7014) ---
$ws.Range("A2").Value = 77
$ws.Range("E2").Value = "SIST_7014"
$ws.Range("F2").Value = "Bases de Datos Geoespaciales"
$ws.Range("G2").Value = "Pucha Cofrep Franz Leonardo"
$ws.Range("H2").Value = "fapucha@utpl.edu.ec"
$ws.Range("I2").Value = 1104483498
$ws.Range("K2").Value = "Ingenierías y Arquitectura"
$ws.Range("L2").Value = "Especialización en Gestión de Geoinformación con mención en Proyectos de Ingeniería"
$ws.Range("M2").Value = "González Jaramillo Víctor Hugo"

# N2 becomes a mailto hyperlink (this shifts the existing V2 hyperlink's
# relationship id down, since N2's is created first)
$ws.Range("V2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("N2"), "mailto:vhgonzalez@utpl.edu.ec", "", "", "mailto:vhgonzalez@utpl.edu.ec")
$ws.Range("N2").Value = "vhgonzalez@utpl.edu.ec"
$ws.Range("N2").Style = "Hipervínculo"

$ws.Range("O2").Value = 1
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = "Especialización"

$ws.Hyperlinks.Add($ws.Range("V2"), "https://utpl.instructure.com/courses/72916")
$ws.Range("V2").Value = "https://utpl.instructure.com/courses/72916"
$ws.Range("V2").Style = "Hipervínculo"

$ws.Range("X2").Value = "SIST_7014_META"
$ws.Range("Y2").Value = "180625 Código banner incorrecto"
$ws.Range("AA2").Value = "SIST_7014"
$ws.Range("AB2").Value = 45841
$ws.Range("AC2").Value = 45841

# --- Sheet view: update selection ---
$ws.Activate()
$ws.Range("R29:R30").Select()
